# Apply the cryptos worksheet refresh (price/volume updates + 3-row coin
# rotation in rows 45-47) as described by the commit's OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; Col=4; Value="64.143.84" },
    @{ Row=2; Col=5; Value="  -0.23%  " },
    @{ Row=3; Col=4; Value="3.322.71" },
    @{ Row=3; Col=5; Value="  -0.78%  " },
    @{ Row=4; Col=5; Value="  +0.14%  " },
    @{ Row=5; Col=4; Value="550.91" },
    @{ Row=5; Col=5; Value="  -0.51%  " },
    @{ Row=6; Col=4; Value="172.47" },
    @{ Row=6; Col=5; Value="  -1.14%  " },
    @{ Row=7; Col=4; Value="0.627" },
    @{ Row=7; Col=5; Value="  +1.98%  " },
    @{ Row=8; Col=5; Value="  +0.02%  " },
    @{ Row=9; Col=4; Value="3.311.38" },
    @{ Row=9; Col=5; Value="  -0.88%  " },
    @{ Row=10; Col=5; Value="  +6.66%  " },
    @{ Row=11; Col=4; Value="0.634" },
    @{ Row=11; Col=5; Value="  +1.25%  " },
    @{ Row=12; Col=4; Value="53.05" },
    @{ Row=12; Col=5; Value="  -1.48%  " },
    @{ Row=13; Col=5; Value="  +1.66%  " },
    @{ Row=14; Col=5; Value="  -0.12%  " },
    @{ Row=15; Col=4; Value="3.855.50" },
    @{ Row=15; Col=5; Value="  -0.56%  " },
    @{ Row=16; Col=5; Value="  +2.39%  " },
    @{ Row=17; Col=5; Value="  -1.67%  " },
    @{ Row=18; Col=4; Value="3.337.23" },
    @{ Row=18; Col=5; Value="  -0.22%  " },
    @{ Row=19; Col=4; Value="64.092.82" },
    @{ Row=19; Col=5; Value="  -0.01%  " },
    @{ Row=20; Col=5; Value="  -1.22%  " },
    @{ Row=21; Col=5; Value="  +0.23%  " },
    @{ Row=22; Col=4; Value="449.18" },
    @{ Row=22; Col=5; Value="  +5.65%  " },
    @{ Row=23; Col=4; Value="5.00" },
    @{ Row=23; Col=5; Value="  +2.66%  " },
    @{ Row=24; Col=4; Value="4.03" },
    @{ Row=24; Col=5; Value="  -1.58%  " },
    @{ Row=25; Col=4; Value="13.95" },
    @{ Row=25; Col=5; Value="  +5.53%  " },
    @{ Row=26; Col=4; Value="86.83" },
    @{ Row=26; Col=5; Value="  +3.40%  " },
    @{ Row=27; Col=5; Value="  +0.97%  " },
    @{ Row=28; Col=4; Value="10.56" },
    @{ Row=28; Col=5; Value="  -1.76%  " },
    @{ Row=29; Col=4; Value="8.56" },
    @{ Row=29; Col=5; Value="  -1.07%  " },
    @{ Row=30; Col=4; Value="30.72" },
    @{ Row=30; Col=5; Value="  +3.33%  " },
    @{ Row=31; Col=4; Value="6.48" },
    @{ Row=31; Col=5; Value="  -2.92%  " },
    @{ Row=32; Col=4; Value="62.67" },
    @{ Row=32; Col=5; Value="  +7.52%  " },
    @{ Row=33; Col=4; Value="11.34" },
    @{ Row=33; Col=5; Value="  -0.76%  " },
    @{ Row=34; Col=4; Value="569.88" },
    @{ Row=34; Col=5; Value="  -0.18%  " },
    @{ Row=35; Col=5; Value="  -1.15%  " },
    @{ Row=36; Col=5; Value="  +0.01%  " },
    @{ Row=37; Col=5; Value="  -0.70%  " },
    @{ Row=38; Col=5; Value="  +0.80%  " },
    @{ Row=39; Col=4; Value="35.12" },
    @{ Row=39; Col=5; Value="  -1.31%  " },
    @{ Row=40; Col=4; Value="0.364" },
    @{ Row=40; Col=5; Value="  -0.59%  " },
    @{ Row=41; Col=5; Value="  -4.05%  " },
    @{ Row=42; Col=4; Value="3.054.53" },
    @{ Row=42; Col=5; Value="  -1.42%  " },
    @{ Row=43; Col=4; Value="0.0412" },
    @{ Row=43; Col=5; Value="  +1.26%  " },
    @{ Row=44; Col=4; Value="2.72" },
    @{ Row=44; Col=5; Value="  -3.30%  " },
    @{ Row=45; Col=2; Value="ApeXProtocol" },
    @{ Row=45; Col=3; Value="https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex" },
    @{ Row=45; Col=4; Value="3.17" },
    @{ Row=45; Col=5; Value="  -1.19%  " },
    @{ Row=46; Col=2; Value="Fetch.AI" },
    @{ Row=46; Col=3; Value="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet" },
    @{ Row=46; Col=4; Value="2.43" },
    @{ Row=46; Col=5; Value="  -1.19%  " },
    @{ Row=47; Col=2; Value="Stellar" },
    @{ Row=47; Col=3; Value="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" },
    @{ Row=47; Col=4; Value="0.133" },
    @{ Row=47; Col=5; Value="  +2.82%  " },
    @{ Row=48; Col=4; Value="1.00" },
    @{ Row=48; Col=5; Value="  +0.16%  " },
    @{ Row=49; Col=4; Value="140.97" },
    @{ Row=49; Col=5; Value="  +3.20%  " },
    @{ Row=50; Col=4; Value="2.52" },
    @{ Row=50; Col=5; Value="  -2.02%  " },
    @{ Row=51; Col=4; Value="8.15" },
    @{ Row=51; Col=5; Value="  -0.84%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    $val = $u.Value
    # Guard against Excel's automatic text->number coercion for strings
    # that look like plain numbers (e.g. "5.00" / "1.00" would otherwise
    # collapse to 5 / 1 and lose the trailing zero). A leading apostrophe
    # forces Excel to keep the literal text, exactly like a human typing
    # '5.00 into the cell would.
    if ($val -match '^-?\d+(\.\d+)?$') {
        $cell.Value = "'" + $val
    } else {
        $cell.Value = $val
    }
}
